$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 25 and row 26
$ws.Range("B25").Value = 6221699
$ws.Range("B26").Value = 6221703
$ws.Range("F25").Value = "FK Maktaaral"
$ws.Range("F26").Value = "Shakhter Karagandy"
$ws.Range("G25").Value = "Kaisar Kyzylorda"
$ws.Range("G26").Value = "FK Aktobe"
$ws.Range("H25").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I25").Value = 2
$ws.Range("I26").Value = 1
$ws.Range("J25").Value = "D"
$ws.Range("J26").Value = "A"
$ws.Range("K25").Value = 3.1
$ws.Range("K26").Value = 3.6
$ws.Range("L25").Value = 3.2
$ws.Range("L26").Value = 3.5
$ws.Range("M25").Value = 2.1
$ws.Range("M26").Value = 1.8
$ws.Range("N25").Value = 2.1
$ws.Range("N26").Value = 3.1
$ws.Range("O25").Value = 3.1
$ws.Range("O26").Value = 3.5
$ws.Range("P25").Value = 3.2
$ws.Range("P26").Value = 1.909
$ws.Range("Q25").Value = -0.25
$ws.Range("Q26").Value = 0.5
$ws.Range("R25").Value = 1.85
$ws.Range("R26").Value = 1.825
$ws.Range("S25").Value = 1.95
$ws.Range("S26").Value = 1.975
$ws.Range("T25").Value = 2.25
$ws.Range("T26").Value = 2.5
$ws.Range("U25").Value = 1.975
$ws.Range("U26").Value = 1.75
$ws.Range("V25").Value = 1.725
$ws.Range("V26").Value = 1.95
$ws.Range("W25").Value = -1
$ws.Range("W26").Value = -1
$ws.Range("X25").Value = 2.1
$ws.Range("X26").Value = -1
$ws.Range("Y25").Value = -1
$ws.Range("Y26").Value = 0.909
$ws.Range("Z25").Value = -0.5
$ws.Range("Z26").Value = -1
$ws.Range("AA25").Value = 0.475
$ws.Range("AA26").Value = 0.9750000000000001
$ws.Range("AB25").Value = 0.9750000000000001
$ws.Range("AB26").Value = -1
$ws.Range("AC25").Value = -1
$ws.Range("AC26").Value = 0.95

# Swap row 37 and row 38
$ws.Range("B37").Value = 6221712
$ws.Range("B38").Value = 6221708
$ws.Range("F37").Value = "FK Aksu"
$ws.Range("F38").Value = "Kaisar Kyzylorda"
$ws.Range("G37").Value = "Shakhter Karagandy"
$ws.Range("G38").Value = "Kairat Almaty"
$ws.Range("H37").Value = 2
$ws.Range("H38").Value = 0
$ws.Range("I37").Value = 1
$ws.Range("I38").Value = 0
$ws.Range("J37").Value = "H"
$ws.Range("J38").Value = "D"
$ws.Range("K37").Value = 2.1
$ws.Range("K38").Value = 3
$ws.Range("L37").Value = 3.25
$ws.Range("L38").Value = 3.4
$ws.Range("M37").Value = 3
$ws.Range("M38").Value = 2.05
$ws.Range("N37").Value = 2.15
$ws.Range("N38").Value = 3.2
$ws.Range("O37").Value = 3.25
$ws.Range("O38").Value = 3.4
$ws.Range("P37").Value = 2.9
$ws.Range("P38").Value = 1.95
$ws.Range("Q37").Value = -0.25
$ws.Range("Q38").Value = 0.5
$ws.Range("R37").Value = 1.95
$ws.Range("R38").Value = 1.75
$ws.Range("S37").Value = 1.85
$ws.Range("S38").Value = 1.95
$ws.Range("T37").Value = 2.5
$ws.Range("T38").Value = 2.25
$ws.Range("U37").Value = 1.975
$ws.Range("U38").Value = 1.925
$ws.Range("V37").Value = 1.825
$ws.Range("V38").Value = 1.875
$ws.Range("W37").Value = 1.15
$ws.Range("W38").Value = -1
$ws.Range("X37").Value = -1
$ws.Range("X38").Value = 2.4
$ws.Range("Y37").Value = -1
$ws.Range("Y38").Value = -1
$ws.Range("Z37").Value = 0.95
$ws.Range("Z38").Value = 0.75
$ws.Range("AA37").Value = -1
$ws.Range("AA38").Value = -1
$ws.Range("AB37").Value = 0.9750000000000001
$ws.Range("AB38").Value = -1
$ws.Range("AC37").Value = -1
$ws.Range("AC38").Value = 0.875

# Swap row 92 and row 93
$ws.Range("B92").Value = 6221748
$ws.Range("B93").Value = 6221749
$ws.Range("F92").Value = "Tobol Kostanay"
$ws.Range("F93").Value = "Kairat Almaty"
$ws.Range("G92").Value = "FK Kaspyi Aktau"
$ws.Range("G93").Value = "FK Aksu"
$ws.Range("H92").Value = 1
$ws.Range("H93").Value = 4
$ws.Range("I92").Value = 1
$ws.Range("I93").Value = 1
$ws.Range("J92").Value = "D"
$ws.Range("J93").Value = "H"
$ws.Range("K92").Value = 1.533
$ws.Range("K93").Value = 1.5
$ws.Range("L92").Value = 3.4
$ws.Range("L93").Value = 3.4
$ws.Range("M92").Value = 6
$ws.Range("M93").Value = 6.5
$ws.Range("N92").Value = 1.444
$ws.Range("N93").Value = 1.333
$ws.Range("O92").Value = 4.2
$ws.Range("O93").Value = 4.5
$ws.Range("P92").Value = 6.5
$ws.Range("P93").Value = 7.5
$ws.Range("Q92").Value = -1.25
$ws.Range("Q93").Value = -1.5
$ws.Range("R92").Value = 1.85
$ws.Range("R93").Value = 1.8
$ws.Range("S92").Value = 1.95
$ws.Range("S93").Value = 2
$ws.Range("T92").Value = 3
$ws.Range("T93").Value = 3
$ws.Range("U92").Value = 1.925
$ws.Range("U93").Value = 1.85
$ws.Range("V92").Value = 1.875
$ws.Range("V93").Value = 1.95
$ws.Range("W92").Value = -1
$ws.Range("W93").Value = 0.333
$ws.Range("X92").Value = 3.2
$ws.Range("X93").Value = -1
$ws.Range("Y92").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z92").Value = -1
$ws.Range("Z93").Value = 0.8
$ws.Range("AA92").Value = 0.95
$ws.Range("AA93").Value = -1
$ws.Range("AB92").Value = -1
$ws.Range("AB93").Value = 0.8500000000000001
$ws.Range("AC92").Value = 0.875
$ws.Range("AC93").Value = -1

# Swap row 99 and row 100
$ws.Range("B99").Value = 6221815
$ws.Range("B100").Value = 6221753
$ws.Range("F99").Value = "FK Atyrau"
$ws.Range("F100").Value = "FK Aksu"
$ws.Range("G99").Value = "Kairat Almaty"
$ws.Range("G100").Value = "Tobol Kostanay"
$ws.Range("H99").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("I100").Value = 3
$ws.Range("J99").Value = "D"
$ws.Range("J100").Value = "A"
$ws.Range("K99").Value = 3
$ws.Range("K100").Value = 2.75
$ws.Range("L99").Value = 3
$ws.Range("L100").Value = 3.1
$ws.Range("M99").Value = 2.25
$ws.Range("M100").Value = 2.375
$ws.Range("N99").Value = 3.1
$ws.Range("N100").Value = 2.625
$ws.Range("O99").Value = 3.1
$ws.Range("O100").Value = 3.2
$ws.Range("P99").Value = 2.15
$ws.Range("P100").Value = 2.45
$ws.Range("Q99").Value = 0.25
$ws.Range("Q100").Value = 0
$ws.Range("R99").Value = 1.85
$ws.Range("R100").Value = 2
$ws.Range("S99").Value = 1.95
$ws.Range("S100").Value = 1.8
$ws.Range("T99").Value = 2.25
$ws.Range("T100").Value = 2.5
$ws.Range("U99").Value = 1.8
$ws.Range("U100").Value = 1.9
$ws.Range("V99").Value = 2
$ws.Range("V100").Value = 1.9
$ws.Range("W99").Value = -1
$ws.Range("W100").Value = -1
$ws.Range("X99").Value = 2.1
$ws.Range("X100").Value = -1
$ws.Range("Y99").Value = -1
$ws.Range("Y100").Value = 1.45
$ws.Range("Z99").Value = 0.425
$ws.Range("Z100").Value = -1
$ws.Range("AA99").Value = -0.5
$ws.Range("AA100").Value = 0.8
$ws.Range("AB99").Value = -1
$ws.Range("AB100").Value = 0.8999999999999999
$ws.Range("AC99").Value = 1
$ws.Range("AC100").Value = -1
